# Apply the "Data retrieved - Mon May 31 19:54:24 UTC 2021" update:
#  - refresh the timestamp on the last existing row (A33)
#  - append a new row (34) with the newly retrieved data point

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the tiny floating-point drift on the existing last row's timestamp.
$ws.Range("A33").Value = 44346.80724647916

# Append the newly retrieved row.
$ws.Range("A34").Value = 44347.82944256221
$ws.Range("B34").Value = 74129
$ws.Range("C34").Value = 62409
$ws.Range("D34").Value = 3218
$ws.Range("E34").Value = 2028
$ws.Range("F34").Value = 1428
$ws.Range("G34").Value = 19305
$ws.Range("H34").Value = 1384
$ws.Range("I34").Value = 817
$ws.Range("J34").Value = 202

# Match the date-formatted style used by the rest of column A.
$ws.Range("A34").NumberFormat = $ws.Range("A33").NumberFormat
